$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 91

# Columns A (Date) and D (Week) contain date/number-looking text that Excel's
# automatic type detection would otherwise coerce into a real date serial /
# number. Force them to be stored as plain text, matching the rest of the
# column, then clear the temporary "@" number format so no stray style is
# left behind on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-29"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "17:38:43"
$ws.Cells.Item($row, 3).Value = "Thursday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "26"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 123283
$ws.Cells.Item($row, 6).Value = 134460
$ws.Cells.Item($row, 7).Value = 163863
$ws.Cells.Item($row, 8).Value = 134147
$ws.Cells.Item($row, 9).Value = 177029
$ws.Cells.Item($row, 10).Value = 115325
$ws.Cells.Item($row, 11).Value = 204635
$ws.Cells.Item($row, 12).Value = 226224
$ws.Cells.Item($row, 13).Value = 176344
$ws.Cells.Item($row, 14).Value = 104510
$ws.Cells.Item($row, 15).Value = 39796
$ws.Cells.Item($row, 16).Value = 33699
$ws.Cells.Item($row, 17).Value = 52571
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36006
$ws.Cells.Item($row, 20).Value = -1
